$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Cells.Item(76, 8).Value = 4014.7144
$ws.Cells.Item(76, 9).Value = 4014.7144
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 4014.7144
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).Value = -3699.7144
$ws.Cells.Item(76, 14).Value = $null
# Row 79
$ws.Cells.Item(79, 8).Value = 4014.7144
$ws.Cells.Item(79, 9).Value = 4014.7144
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 4014.7144
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).Value = -2922.7144
$ws.Cells.Item(79, 14).Value = $null
# Row 132
$ws.Cells.Item(132, 8).Value = 11114351
$ws.Cells.Item(132, 9).Value = 15876588
$ws.Cells.Item(132, 10).Value = 2464
$ws.Cells.Item(132, 11).Value = 47629764
$ws.Cells.Item(132, 12).Value = 7392
$ws.Cells.Item(132, 13).Value = -47627234
$ws.Cells.Item(132, 14).Value = -12452
# Row 137
$ws.Cells.Item(137, 8).Value = 2328.577
$ws.Cells.Item(137, 9).Value = 1358.7778
$ws.Cells.Item(137, 11).Value = 4076.3334
$ws.Cells.Item(137, 13).Value = -1526.3334
# Row 138
$ws.Cells.Item(138, 8).Value = 1744.763
$ws.Cells.Item(138, 9).Value = 840.36365
$ws.Cells.Item(138, 10).Value = 1860.4419
$ws.Cells.Item(138, 11).Value = 2521.09095
$ws.Cells.Item(138, 12).Value = 5581.3257
$ws.Cells.Item(138, 13).Value = 2618.90905
$ws.Cells.Item(138, 14).Value = -15861.3257
# Row 140
$ws.Cells.Item(140, 8).Value = 34853.332
$ws.Cells.Item(140, 9).Value = 3000
$ws.Cells.Item(140, 10).Value = 50780
$ws.Cells.Item(140, 11).Value = 3000
$ws.Cells.Item(140, 12).Value = 50780
$ws.Cells.Item(140, 13).Value = 2180
$ws.Cells.Item(140, 14).Value = -61140
# Row 141
$ws.Cells.Item(141, 8).Value = 2220.7144
$ws.Cells.Item(141, 9).Value = 2220.7144
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 6662.1432
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = -1482.1432
$ws.Cells.Item(141, 14).Value = $null

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Cells.Item(4, 8).Value = 92.85714
$ws.Cells.Item(4, 9).Value = 87.5
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 87.5
$ws.Cells.Item(4, 12).Value = 100
$ws.Cells.Item(4, 13).Value = 28.5
$ws.Cells.Item(4, 14).Value = -332
# Row 32
$ws.Cells.Item(32, 8).Value = 2689.8572
$ws.Cells.Item(32, 9).Value = 2873.4092
$ws.Cells.Item(32, 11).Value = 2873.4092
$ws.Cells.Item(32, 13).Value = -2586.4092
# Row 61
$ws.Cells.Item(61, 8).Value = 1608.75
$ws.Cells.Item(61, 9).Value = 1608.75
$ws.Cells.Item(61, 11).Value = 1608.75
$ws.Cells.Item(61, 13).Value = -1396.75
# Row 82
$ws.Cells.Item(82, 8).Value = 27863.5
$ws.Cells.Item(82, 10).Value = 27863.5
$ws.Cells.Item(82, 12).Value = 27863.5
$ws.Cells.Item(82, 14).Value = -28585.5
# Row 85
$ws.Cells.Item(85, 8).Value = 27863.5
$ws.Cells.Item(85, 10).Value = 27863.5
$ws.Cells.Item(85, 12).Value = 27863.5
$ws.Cells.Item(85, 14).Value = -30359.5
# Row 102
$ws.Cells.Item(102, 8).Value = 13897688
$ws.Cells.Item(102, 9).Value = 15161041
$ws.Cells.Item(102, 11).Value = 15161041
$ws.Cells.Item(102, 13).Value = -15159419
# Row 114
$ws.Cells.Item(114, 8).Value = 24199.2
$ws.Cells.Item(114, 10).Value = 24199.2
$ws.Cells.Item(114, 12).Value = 24199.2
$ws.Cells.Item(114, 14).Value = -32877.2
# Row 136
$ws.Cells.Item(136, 8).Value = 1608.75
$ws.Cells.Item(136, 9).Value = 1608.75
$ws.Cells.Item(136, 11).Value = 4826.25
$ws.Cells.Item(136, 13).Value = -2276.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 7353525.5
$ws.Cells.Item(94, 9).Value = 10870043
$ws.Cells.Item(94, 10).Value = 808.0909
$ws.Cells.Item(94, 11).Value = 10870043
$ws.Cells.Item(94, 12).Value = 808.0909
$ws.Cells.Item(94, 13).Value = -10869592
$ws.Cells.Item(94, 14).Value = -1710.0909

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1497.4783
$ws.Cells.Item(31, 9).Value = 1889.2
$ws.Cells.Item(31, 10).Value = 1388.6666
$ws.Cells.Item(31, 11).Value = 1889.2
$ws.Cells.Item(31, 12).Value = 1388.6666
$ws.Cells.Item(31, 13).Value = -1594.2
$ws.Cells.Item(31, 14).Value = -1978.6666
# Row 34
$ws.Cells.Item(34, 8).Value = 1497.4783
$ws.Cells.Item(34, 9).Value = 1889.2
$ws.Cells.Item(34, 10).Value = 1388.6666
$ws.Cells.Item(34, 11).Value = 1889.2
$ws.Cells.Item(34, 12).Value = 1388.6666
$ws.Cells.Item(34, 13).Value = -1687.2
$ws.Cells.Item(34, 14).Value = -1792.6666
# Row 99
$ws.Cells.Item(99, 8).Value = 1756115.6
$ws.Cells.Item(99, 9).Value = 2925570.5
$ws.Cells.Item(99, 10).Value = 1933.3334
$ws.Cells.Item(99, 11).Value = 2925570.5
$ws.Cells.Item(99, 12).Value = 1933.3334
$ws.Cells.Item(99, 13).Value = -2924072.5
$ws.Cells.Item(99, 14).Value = -4929.3334
# Row 123
$ws.Cells.Item(123, 8).Value = 79889.625
$ws.Cells.Item(123, 10).Value = 79889.625
$ws.Cells.Item(123, 12).Value = 79889.625
$ws.Cells.Item(123, 14).Value = -89689.625
# Row 126
$ws.Cells.Item(126, 8).Value = 1756115.6
$ws.Cells.Item(126, 9).Value = 2925570.5
$ws.Cells.Item(126, 10).Value = 1933.3334
$ws.Cells.Item(126, 11).Value = 8776711.5
$ws.Cells.Item(126, 12).Value = 5800.0002
$ws.Cells.Item(126, 13).Value = -8774241.5
$ws.Cells.Item(126, 14).Value = -10740.0002
# Row 132
$ws.Cells.Item(132, 8).Value = 5102.1562
$ws.Cells.Item(132, 9).Value = 5432.9614
$ws.Cells.Item(132, 10).Value = 3668.6667
$ws.Cells.Item(132, 11).Value = 16298.8842
$ws.Cells.Item(132, 12).Value = 11006.0001
$ws.Cells.Item(132, 13).Value = -13768.8842
$ws.Cells.Item(132, 14).Value = -16066.0001
# Row 134
$ws.Cells.Item(134, 8).Value = 11495192
$ws.Cells.Item(134, 9).Value = 11495192
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 34485576
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -34483041
$ws.Cells.Item(134, 14).Value = $null

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Cells.Item(68, 8).Value = 1729.975
$ws.Cells.Item(68, 9).Value = 687.75
$ws.Cells.Item(68, 10).Value = 1990.5312
$ws.Cells.Item(68, 11).Value = 2063.25
$ws.Cells.Item(68, 12).Value = 5971.5936
$ws.Cells.Item(68, 13).Value = -1252.25
$ws.Cells.Item(68, 14).Value = -7593.5936
# Row 71
$ws.Cells.Item(71, 8).Value = 1729.975
$ws.Cells.Item(71, 9).Value = 687.75
$ws.Cells.Item(71, 10).Value = 1990.5312
$ws.Cells.Item(71, 11).Value = 6189.75
$ws.Cells.Item(71, 12).Value = 17914.7808
$ws.Cells.Item(71, 13).Value = -2133.75
$ws.Cells.Item(71, 14).Value = -26026.7808
# Row 121
$ws.Cells.Item(121, 8).Value = 510.33334
$ws.Cells.Item(121, 9).Value = 254.33333
$ws.Cells.Item(121, 10).Value = 766.3333
$ws.Cells.Item(121, 11).Value = 762.99999
$ws.Cells.Item(121, 12).Value = 2298.9999
$ws.Cells.Item(121, 13).Value = 547.00001
$ws.Cells.Item(121, 14).Value = -4918.9999
# Row 131
$ws.Cells.Item(131, 8).Value = 20001420
$ws.Cells.Item(131, 10).Value = 1606.5
$ws.Cells.Item(131, 12).Value = 4819.5
$ws.Cells.Item(131, 14).Value = -14899.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Cells.Item(5, 8).Value = 2000
$ws.Cells.Item(5, 9).Value = 2000
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 2000
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -1888
$ws.Cells.Item(5, 14).Value = $null
# Row 111
$ws.Cells.Item(111, 8).Value = 100000000
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 14).Value = $null
# Row 122
$ws.Cells.Item(122, 8).Value = 3144.3125
$ws.Cells.Item(122, 9).Value = 2030.9
$ws.Cells.Item(122, 11).Value = 6092.700000000001
$ws.Cells.Item(122, 13).Value = -3642.700000000001
# Row 132
$ws.Cells.Item(132, 8).Value = 4019.7
$ws.Cells.Item(132, 9).Value = 3671.2856
$ws.Cells.Item(132, 11).Value = 11013.8568
$ws.Cells.Item(132, 13).Value = -8483.856800000001
# Row 140
$ws.Cells.Item(140, 8).Value = 70000
$ws.Cells.Item(140, 10).Value = 70000
$ws.Cells.Item(140, 12).Value = 70000
$ws.Cells.Item(140, 14).Value = -80360
# Row 141
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).Value = $null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Cells.Item(16, 8).Value = 1597.2858
$ws.Cells.Item(16, 9).Value = 1405.6364
$ws.Cells.Item(16, 10).Value = 2300
$ws.Cells.Item(16, 11).Value = 1405.6364
$ws.Cells.Item(16, 12).Value = 2300
$ws.Cells.Item(16, 13).Value = -1235.6364
$ws.Cells.Item(16, 14).Value = -2640
# Row 100
$ws.Cells.Item(100, 8).Value = 1919.1818
$ws.Cells.Item(100, 9).Value = 1867.6666
$ws.Cells.Item(100, 11).Value = 1867.6666
$ws.Cells.Item(100, 13).Value = -1326.6666
# Row 136
$ws.Cells.Item(136, 8).Value = 7867.1333
$ws.Cells.Item(136, 9).Value = 10999.4
$ws.Cells.Item(136, 10).Value = 1602.6
$ws.Cells.Item(136, 11).Value = 32998.2
$ws.Cells.Item(136, 12).Value = 4807.799999999999
$ws.Cells.Item(136, 13).Value = -30448.2
$ws.Cells.Item(136, 14).Value = -9907.799999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 41672200
$ws.Cells.Item(62, 9).Value = 50004260
$ws.Cells.Item(62, 11).Value = 50004260
$ws.Cells.Item(62, 13).Value = -50003636
# Row 65
$ws.Cells.Item(65, 8).Value = 41672200
$ws.Cells.Item(65, 9).Value = 50004260
$ws.Cells.Item(65, 11).Value = 250021300
$ws.Cells.Item(65, 13).Value = -250018180
# Row 100
$ws.Cells.Item(100, 8).Value = 1482.5
$ws.Cells.Item(100, 9).Value = 971
$ws.Cells.Item(100, 10).Value = 2249.75
$ws.Cells.Item(100, 11).Value = 1942
$ws.Cells.Item(100, 12).Value = 4499.5
$ws.Cells.Item(100, 13).Value = -1401
$ws.Cells.Item(100, 14).Value = -5581.5
# Row 107
$ws.Cells.Item(107, 8).Value = 867.5
$ws.Cells.Item(107, 9).Value = 834
$ws.Cells.Item(107, 11).Value = 2502
$ws.Cells.Item(107, 13).Value = -582
# Row 122
$ws.Cells.Item(122, 8).Value = 32509388
$ws.Cells.Item(122, 9).Value = 37153300
$ws.Cells.Item(122, 11).Value = 111459900
$ws.Cells.Item(122, 13).Value = -111457450
# Row 126
$ws.Cells.Item(126, 8).Value = 35843390
$ws.Cells.Item(126, 9).Value = 61729070
$ws.Cells.Item(126, 11).Value = 185187210
$ws.Cells.Item(126, 13).Value = -185184740
# Row 132
$ws.Cells.Item(132, 8).Value = 3018.4443
$ws.Cells.Item(132, 9).Value = 2608.1
$ws.Cells.Item(132, 10).Value = 4190.857
$ws.Cells.Item(132, 11).Value = 7824.299999999999
$ws.Cells.Item(132, 12).Value = 12572.571
$ws.Cells.Item(132, 13).Value = -5294.299999999999
$ws.Cells.Item(132, 14).Value = -17632.571
# Row 136
$ws.Cells.Item(136, 8).Value = 882.13635
$ws.Cells.Item(136, 9).Value = 609.5
$ws.Cells.Item(136, 10).Value = 984.375
$ws.Cells.Item(136, 11).Value = 1828.5
$ws.Cells.Item(136, 12).Value = 2953.125
$ws.Cells.Item(136, 13).Value = 721.5
$ws.Cells.Item(136, 14).Value = -8053.125
# Row 140
$ws.Cells.Item(140, 8).Value = 31858
$ws.Cells.Item(140, 10).Value = 31858
$ws.Cells.Item(140, 12).Value = 31858
$ws.Cells.Item(140, 14).Value = -42218
# Row 141
$ws.Cells.Item(141, 8).Value = 44700
$ws.Cells.Item(141, 10).Value = 44700
$ws.Cells.Item(141, 12).Value = 44700
$ws.Cells.Item(141, 14).Value = -55060
